$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue 'D2' '57.209.87'
Set-TextValue 'E2' '  +0.05%  '
Set-TextValue 'D3' '2.406.19'
Set-TextValue 'E3' '  -4.06%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '488.47'
Set-TextValue 'E5' '  -1.16%  '
Set-TextValue 'D6' '154.17'
Set-TextValue 'E6' '  +0.12%  '
Set-TextValue 'D7' '0.615'
Set-TextValue 'E7' '  +19.12%  '
Set-TextValue 'D8' '0.998'
Set-TextValue 'E8' '  +0.19%  '
Set-TextValue 'D9' '2.422.79'
Set-TextValue 'E9' '  -3.99%  '
Set-TextValue 'D10' '6.28'
Set-TextValue 'E10' '  +8.80%  '
Set-TextValue 'D11' '0.0996'
Set-TextValue 'E11' '  -0.17%  '
Set-TextValue 'D12' '0.333'
Set-TextValue 'E12' '  -1.40%  '
Set-TextValue 'E13' '  +1.27%  '
Set-TextValue 'D14' '2.829.66'
Set-TextValue 'E14' '  -3.72%  '
Set-TextValue 'D15' '57.178.80'
Set-TextValue 'E15' '  -0.19%  '
Set-TextValue 'D16' '20.56'
Set-TextValue 'E16' '  -3.56%  '
Set-TextValue 'E17' '  -3.53%  '
Set-TextValue 'D18' '2.423.17'
Set-TextValue 'E18' '  -4.09%  '
Set-TextValue 'D19' '4.71'
Set-TextValue 'E19' '  +3.35%  '
Set-TextValue 'D20' '324.44'
Set-TextValue 'E20' '  +0.35%  '
Set-TextValue 'D21' '9.95'
Set-TextValue 'E21' '  -3.81%  '
Set-TextValue 'D22' '0.996'
Set-TextValue 'E22' '  -0.16%  '
Set-TextValue 'D23' '5.91'
Set-TextValue 'E23' '  -0.19%  '
Set-TextValue 'D24' '57.78'
Set-TextValue 'E24' '  -1.30%  '
Set-TextValue 'D25' '0.404'
Set-TextValue 'E25' '  -1.69%  '
Set-TextValue 'D26' '0.997'
Set-TextValue 'E26' '  -0.19%  '
Set-TextValue 'E27' '  -1.63%  '
Set-TextValue 'D28' '2.514.83'
Set-TextValue 'E28' '  -3.71%  '
Set-TextValue 'E29' '  -5.54%  '
Set-TextValue 'D30' '0.0₃0779'
Set-TextValue 'E30' '  -6.01%  '
Set-TextValue 'D31' '1.00'
Set-TextValue 'E31' '  +0.10%  '
Set-TextValue 'D32' '150.39'
Set-TextValue 'E32' '  -0.98%  '
Set-TextValue 'D33' '18.54'
Set-TextValue 'E33' '  +1.25%  '
Set-TextValue 'E34' '  -0.69%  '
Set-TextValue 'E35' '  -0.61%  '
Set-TextValue 'D36' '1.15'
Set-TextValue 'E36' '  -1.29%  '
Set-TextValue 'D37' '3.75'
Set-TextValue 'E37' '  -2.24%  '
Set-TextValue 'D38' '0.837'
Set-TextValue 'E38' '  -5.59%  '
Set-TextValue 'E39' '  +8.45%  '
Set-TextValue 'E40' '  -0.63%  '
Set-TextValue 'D41' '3.52'
Set-TextValue 'E41' '  -0.44%  '
Set-TextValue 'E42' '  -2.81%  '
Set-TextValue 'D43' '0.997'
Set-TextValue 'E43' '  +0.35%  '
Set-TextValue 'D44' '273.08'
Set-TextValue 'E44' '  +1.39%  '
Set-TextValue 'E45' '  -3.97%  '
Set-TextValue 'D46' '0.0529'
Set-TextValue 'E46' '  -6.24%  '
Set-TextValue 'D47' '10.21'
Set-TextValue 'E47' '  +0.05%  '
Set-TextValue 'D48' '0.0227'
Set-TextValue 'E48' '  -1.37%  '
Set-TextValue 'D49' '4.48'
Set-TextValue 'E49' '  -9.05%  '
Set-TextValue 'D50' '1.883.88'
Set-TextValue 'E50' '  -0.95%  '
Set-TextValue 'D51' '17.45'
Set-TextValue 'E51' '  -3.34%  '
